$d = $word.ActiveDocument

# --- Step 1: Remove the "Meta description: ..." paragraph that currently
# follows the H1 title at the top of the document. ---
$metaPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^Meta description:") {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# --- Step 2: Locate the "Limited number of pay lines" bullet (the paragraph
# that immediately precedes the trailing DALLE image-prompt paragraph) and
# append a brand-new bold paragraph right after it - i.e. right before the
# DALLE paragraph. Inserting the text + paragraph mark directly at the end
# of this plain (non-italic, non-bold) paragraph means the new paragraph
# naturally picks up plain formatting, instead of inheriting the italic
# formatting of the DALLE paragraph that follows. ---
$limitedPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^Limited number of pay lines") {
        $limitedPara = $p
        break
    }
}

if ($limitedPara -ne $null) {
    $endRange = $limitedPara.Range.Duplicate
    $endRange.Collapse(0)  # wdCollapseEnd
    $cr = [char]13
    $endRange.InsertAfter("Play Bananas Go Bahamas for Free - A Festive Slot Game" + $cr)

    $newPara = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq ("Play Bananas Go Bahamas for Free - A Festive Slot Game" + $cr)) {
            $newPara = $p
        }
    }

    # Bold only the text itself, not the trailing paragraph mark, so the
    # paragraph-mark run properties stay untouched (matching plain text
    # paragraphs elsewhere in the document).
    $textOnly = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
    $textOnly.Font.Bold = $true
}

# Replace the old DALLE prompt text with the new meta-description text,
# keeping the existing italic formatting of that paragraph intact.
$d.Content.Find.Execute(
    "Create a Cartoon Image of Happy Maya Warrior with Glasses for Bananas Go Bahamas DALLE, we need you to create a feature image for the slot game " + [char]34 + "Bananas Go Bahamas" + [char]34 + " that is in cartoon style and features a happy Maya warrior with glasses. The image should be colorful, vibrant and have a tropical feel to it. The Maya warrior should be shown as happy and carefree, holding a banana in one hand and wearing sunglasses. The background should feature palm trees, sandy beaches, and crystal clear waters. The image should be eye-catching and represent the fun, tropical vibe of the slot game. Remember to keep the image lighthearted and playful, in keeping with the other characters in the game. Good luck!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Enjoy the festive atmosphere of Bananas Go Bahamas with its varied fruit symbols and opportunities to win free spins.",
    2
)
